$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '262.18'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.70%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.75'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.79%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.693'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.17%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06096'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.03%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.701'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.45%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8503'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.08%'

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9090'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.66%'

# Row 9
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1407'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.05%'

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05063'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.92%'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07102'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.25%'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03149'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.41%'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09050'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.03%'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001533'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.23%'

# Row 15
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006186'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.04%'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006026'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.80%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.449'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.07%'

# Row 18
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.164'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.18%'

# Row 19
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.167'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.61%'

# Row 20
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3072'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.19%'

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.32%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.084'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.13%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04240'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001177'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-3.17%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004058'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '6.89%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03941'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1112'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.08%'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.004180'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2.30%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002110'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.94%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01149'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-29.63%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005102'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.35%'

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.00%'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2576'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '55.71%'

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.00%'

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.00%'
